$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$curlyApos = [char]0x2019

# Update the lecture topics/content for the latter part of the schedule.
# (Row order matters for shared-string append ordering, so write top to bottom.)
$ws.Range("B34").Value = "Doing stuff with Arrays in a loop & Nested Loops!   | Arrays Part2 HW && Arrays Part3  | Go over Arrays Part1"
$ws.Range("B35").Value = " Lists!  Life just got easier.  |  Lists HW |  Go over Arrays Part2 & 3"
$ws.Range("B37").Value = "Crash-cource on OOP"
$ws.Range("B39").Value = "IO -> let" + $curlyApos + "s mess with some text files!  & Try/catch  (HW -> Load files and print duplicate numbers)"
$ws.Range("B45").Value = "Learning GUI stuff -> GUI HW "
$ws.Range("B47").Value = "Stacks and Queues   ->  HW  ->  create a Queing system  (using a GUI)"
$ws.Range("B55").Value = ""

# Row 34 grew taller to fit the longer, wrapped text.
$ws.Rows("34").RowHeight = 39.75

# Move the active selection to the new edit location (matches saved view state).
$ws.Range("B51:B52").Select()
